# Update cryptocurrency price/volume data to the latest scraped values.
# (GitHub Actions scheduled refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.423.80'
$ws.Range("E2").Value = '  +0.30%  '
# Row 3
$ws.Range("D3").Value = '1.868.34'
$ws.Range("E3").Value = '  -0.50%  '
# Row 4
$ws.Range("E4").Value = '  -0.14%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.85'
$ws.Range("E5").Value = '  +0.63%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7073'
$ws.Range("E6").Value = '  -0.47%  '
# Row 7
$ws.Range("E7").Value = '  -0.15%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07879'
$ws.Range("E8").Value = '  -1.45%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3133'
$ws.Range("E9").Value = '  -0.80%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.55'
$ws.Range("E10").Value = '  -1.58%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07960'
$ws.Range("E11").Value = '  -3.97%  '
# Row 12
$ws.Range("D12").Value = '1.889.80'
$ws.Range("E12").Value = '  +0.10%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.215'
$ws.Range("E13").Value = '  -0.65%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.43'
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7017'
$ws.Range("E15").Value = '  -1.69%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.521'
$ws.Range("E16").Value = '  +2.49%  '
# Row 17
$ws.Range("D17").Value = '29.446.26'
$ws.Range("E17").Value = '  +0.31%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008357'
$ws.Range("E18").Value = '  -1.88%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.40'
$ws.Range("E19").Value = '  +3.45%  '
# Row 20
$ws.Range("D20").Value = '2.125.24'
$ws.Range("E20").Value = '  -0.73%  '
# Row 21
$ws.Range("E21").Value = '  -1.06%  '
# Row 22
$ws.Range("E22").Value = '  -0.13%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.634'
$ws.Range("E23").Value = '  -1.88%  '
# Row 24
$ws.Range("E24").Value = '  -0.17%  '
# Row 25
$ws.Range("E25").Value = '  -0.15%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.006'
$ws.Range("E26").Value = '  -0.70%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.48'
$ws.Range("E27").Value = '  -0.74%  '
# Row 28
$ws.Range("E28").Value = '  +0.99%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.500'
$ws.Range("E29").Value = '  -0.28%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.333'
$ws.Range("E30").Value = '  -1.88%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.256'
$ws.Range("E31").Value = '  -1.54%  '
# Row 32
$ws.Range("E32").Value = '  +1.32%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05318'
$ws.Range("E33").Value = '  -1.02%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.894'
$ws.Range("E34").Value = '  -2.10%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7495'
$ws.Range("E35").Value = '  -1.93%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.174'
$ws.Range("E36").Value = '  -0.90%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  +0.98%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01891'
$ws.Range("E38").Value = '  +0.41%  '
# Row 39
$ws.Range("D39").Value = '1.273.93'
$ws.Range("E39").Value = '  +1.05%  '
# Row 40
$ws.Range("E40").Value = '  -0.05%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8947'
$ws.Range("E41").Value = '  -1.09%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.066'
$ws.Range("E42").Value = '  -6.86%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '109.13'
$ws.Range("E43").Value = '  -3.40%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.34'
$ws.Range("E44").Value = '  -3.94%  '
# Row 45
$ws.Range("E45").Value = '  -0.17%  '
# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.030.68'
$ws.Range("E46").Value = '  +0.04%  '
# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000126'
$ws.Range("E47").Value = '  -4.71%  '
# Row 48
$ws.Range("E48").Value = '  -0.10%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.562'
$ws.Range("E49").Value = '  +1.21%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5183'
$ws.Range("E50").Value = '  -0.90%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4307'
$ws.Range("E51").Value = '  -1.49%  '
